# fix: pre_post IAHP addon features assignment
#
# Updates the settings sheet:
#   - project_name:  iahp_panel_2 -> iahp_panel_1
#   - data_subsets:  Granulos -> B, CD4_T, CD8_T, Monos_and_DCs, NK, TCRgd_T
#   - umap_n:        15 -> 10
#   - umap_min_dist: 0.1 -> 0.2
# and moves the active selection from B22 to B38 to match where the user
# ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "iahp_panel_1"
$ws.Range("B22").Value = "B, CD4_T, CD8_T, Monos_and_DCs, NK, TCRgd_T"
$ws.Range("B36").Value = 10
$ws.Range("B37").Value = 0.2

# Reflect the final cursor position / selection recorded in the workbook.
$ws.Range("B38").Select()
